$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GENERAL_INTENTS")
$ws.Activate()

$ws.Range("A2").Value = "I want to transfer AMOUNT to PERSON in BANK_ACC bank account."
$ws.Range("A3").Value = "I want to pay bill of AMOUNT with PERSON in BANK_ACC bank account"
$ws.Range("A4").Value = "I want to deposit AMOUNT to BANK_ACC bank account"
$ws.Range("A5").Value = "I want to check balance in BANK_ACC bank account"

$ws.Range("C4").Select()
